# This script reproduces a fresh scrape update for the "LÍNEA 141" schedule
# workbook (commit: "Horarios actualizados Línea 141 - 511").
# It updates the header metadata (last-updated timestamp / row count) and
# rewrites each sheet's data block (sorted by Hora_Llegada), including the
# newly scraped rows appended at the end of sheets 1 and 2.

$wb = $excel.ActiveWorkbook

# ---- Worksheet: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 12:32:47"
$ws1.Range("A3").Value = "Total filas: 239"
$ws1_data = New-Object "object[,]" 239,5
$ws1_data[0,0] = "05:42:22"; $ws1_data[0,1] = "05:43"; $ws1_data[0,2] = "14_ABASTO"; $ws1_data[0,3] = 1; $ws1_data[0,4] = "LP1912"
$ws1_data[1,0] = "05:42:22"; $ws1_data[1,1] = "05:52"; $ws1_data[1,2] = "17_ROMERO"; $ws1_data[1,3] = 10; $ws1_data[1,4] = "LP1912"
$ws1_data[2,0] = "05:42:22"; $ws1_data[2,1] = "06:01"; $ws1_data[2,2] = "16_SANTA ANA"; $ws1_data[2,3] = 19; $ws1_data[2,4] = "LP1912"
$ws1_data[3,0] = "05:42:22"; $ws1_data[3,1] = "06:04"; $ws1_data[3,2] = "10_OLMOS"; $ws1_data[3,3] = 22; $ws1_data[3,4] = "LP1912"
$ws1_data[4,0] = "05:42:22"; $ws1_data[4,1] = "06:11"; $ws1_data[4,2] = "215A_EL PATO"; $ws1_data[4,3] = 29; $ws1_data[4,4] = "LP1912"
$ws1_data[5,0] = "05:42:22"; $ws1_data[5,1] = "06:24"; $ws1_data[5,2] = "11_ETCHEVERRY"; $ws1_data[5,3] = 42; $ws1_data[5,4] = "LP1912"
$ws1_data[6,0] = "05:42:22"; $ws1_data[6,1] = "06:27"; $ws1_data[6,2] = "23_HERNANDEZ"; $ws1_data[6,3] = 45; $ws1_data[6,4] = "LP1912"
$ws1_data[7,0] = "05:42:22"; $ws1_data[7,1] = "06:31"; $ws1_data[7,2] = "16_SANTA ANA"; $ws1_data[7,3] = 49; $ws1_data[7,4] = "LP1912"
$ws1_data[8,0] = "05:42:22"; $ws1_data[8,1] = "06:31"; $ws1_data[8,2] = "17X38_ROMERO"; $ws1_data[8,3] = 49; $ws1_data[8,4] = "LP1912"
$ws1_data[9,0] = "05:42:22"; $ws1_data[9,1] = "06:39"; $ws1_data[9,2] = "225_C ROCA-H SUR"; $ws1_data[9,3] = 57; $ws1_data[9,4] = "LP1912"
$ws1_data[10,0] = "06:33:46"; $ws1_data[10,1] = "06:50"; $ws1_data[10,2] = "215A_EL PATO"; $ws1_data[10,3] = 17; $ws1_data[10,4] = "LP1912"
$ws1_data[11,0] = "06:16:15"; $ws1_data[11,1] = "06:51"; $ws1_data[11,2] = "215A_EL PATO"; $ws1_data[11,3] = 35; $ws1_data[11,4] = "LP1912"
$ws1_data[12,0] = "06:52:52"; $ws1_data[12,1] = "06:52"; $ws1_data[12,2] = "215A_EL PATO"; $ws1_data[12,3] = 0; $ws1_data[12,4] = "LP1912"
$ws1_data[13,0] = "05:42:22"; $ws1_data[13,1] = "06:54"; $ws1_data[13,2] = "14_ABASTO"; $ws1_data[13,3] = 72; $ws1_data[13,4] = "LP1912"
$ws1_data[14,0] = "06:52:52"; $ws1_data[14,1] = "06:55"; $ws1_data[14,2] = "14_ABASTO"; $ws1_data[14,3] = 3; $ws1_data[14,4] = "LP1912"
$ws1_data[15,0] = "06:45:50"; $ws1_data[15,1] = "07:00"; $ws1_data[15,2] = "16_SANTA ANA"; $ws1_data[15,3] = 15; $ws1_data[15,4] = "LP1912"
$ws1_data[16,0] = "05:57:08"; $ws1_data[16,1] = "07:01"; $ws1_data[16,2] = "16_SANTA ANA"; $ws1_data[16,3] = 64; $ws1_data[16,4] = "LP1912"
$ws1_data[17,0] = "06:45:50"; $ws1_data[17,1] = "07:03"; $ws1_data[17,2] = "225_GOMEZ"; $ws1_data[17,3] = 18; $ws1_data[17,4] = "LP1912"
$ws1_data[18,0] = "05:42:22"; $ws1_data[18,1] = "07:04"; $ws1_data[18,2] = "225_GOMEZ"; $ws1_data[18,3] = 82; $ws1_data[18,4] = "LP1912"
$ws1_data[19,0] = "06:33:46"; $ws1_data[19,1] = "07:06"; $ws1_data[19,2] = "215C_EL PATO"; $ws1_data[19,3] = 33; $ws1_data[19,4] = "LP1912"
$ws1_data[20,0] = "05:42:22"; $ws1_data[20,1] = "07:07"; $ws1_data[20,2] = "215C_EL PATO"; $ws1_data[20,3] = 85; $ws1_data[20,4] = "LP1912"
$ws1_data[21,0] = "06:33:46"; $ws1_data[21,1] = "07:13"; $ws1_data[21,2] = "14X44_ABASTO"; $ws1_data[21,3] = 40; $ws1_data[21,4] = "LP1912"
$ws1_data[22,0] = "05:42:22"; $ws1_data[22,1] = "07:14"; $ws1_data[22,2] = "14X44_ABASTO"; $ws1_data[22,3] = 92; $ws1_data[22,4] = "LP1912"
$ws1_data[23,0] = "06:33:46"; $ws1_data[23,1] = "07:15"; $ws1_data[23,2] = "16_SANTA ANA"; $ws1_data[23,3] = 42; $ws1_data[23,4] = "LP1912"
$ws1_data[24,0] = "06:33:46"; $ws1_data[24,1] = "07:20"; $ws1_data[24,2] = "215A_EL PATO"; $ws1_data[24,3] = 47; $ws1_data[24,4] = "LP1912"
$ws1_data[25,0] = "05:42:22"; $ws1_data[25,1] = "07:21"; $ws1_data[25,2] = "215A_EL PATO"; $ws1_data[25,3] = 99; $ws1_data[25,4] = "LP1912"
$ws1_data[26,0] = "06:52:52"; $ws1_data[26,1] = "07:28"; $ws1_data[26,2] = "14_ABASTO"; $ws1_data[26,3] = 36; $ws1_data[26,4] = "LP1912"
$ws1_data[27,0] = "05:57:08"; $ws1_data[27,1] = "07:29"; $ws1_data[27,2] = "14_ABASTO"; $ws1_data[27,3] = 92; $ws1_data[27,4] = "LP1912"
$ws1_data[28,0] = "05:42:22"; $ws1_data[28,1] = "07:33"; $ws1_data[28,2] = "23_HERNANDEZ"; $ws1_data[28,3] = 111; $ws1_data[28,4] = "LP1912"
$ws1_data[29,0] = "05:57:08"; $ws1_data[29,1] = "07:34"; $ws1_data[29,2] = "23_HERNANDEZ"; $ws1_data[29,3] = 97; $ws1_data[29,4] = "LP1912"
$ws1_data[30,0] = "06:45:50"; $ws1_data[30,1] = "07:35"; $ws1_data[30,2] = "17X38_ROMERO"; $ws1_data[30,3] = 50; $ws1_data[30,4] = "LP1912"
$ws1_data[31,0] = "05:42:22"; $ws1_data[31,1] = "07:36"; $ws1_data[31,2] = "17X38_ROMERO"; $ws1_data[31,3] = 114; $ws1_data[31,4] = "LP1912"
$ws1_data[32,0] = "06:33:46"; $ws1_data[32,1] = "07:36"; $ws1_data[32,2] = "27_EL RETIRO"; $ws1_data[32,3] = 63; $ws1_data[32,4] = "LP1912"
$ws1_data[33,0] = "05:42:22"; $ws1_data[33,1] = "07:37"; $ws1_data[33,2] = "27_EL RETIRO"; $ws1_data[33,3] = 115; $ws1_data[33,4] = "LP1912"
$ws1_data[34,0] = "07:36:59"; $ws1_data[34,1] = "07:38"; $ws1_data[34,2] = "16_SANTA ANA"; $ws1_data[34,3] = 2; $ws1_data[34,4] = "LP1912"
$ws1_data[35,0] = "06:33:46"; $ws1_data[35,1] = "07:43"; $ws1_data[35,2] = "10_OLMOS"; $ws1_data[35,3] = 70; $ws1_data[35,4] = "LP1912"
$ws1_data[36,0] = "05:57:08"; $ws1_data[36,1] = "07:44"; $ws1_data[36,2] = "10_OLMOS"; $ws1_data[36,3] = 107; $ws1_data[36,4] = "LP1912"
$ws1_data[37,0] = "07:36:59"; $ws1_data[37,1] = "07:50"; $ws1_data[37,2] = "15_ABASTO"; $ws1_data[37,3] = 14; $ws1_data[37,4] = "LP1912"
$ws1_data[38,0] = "05:57:08"; $ws1_data[38,1] = "07:51"; $ws1_data[38,2] = "15_ABASTO"; $ws1_data[38,3] = 114; $ws1_data[38,4] = "LP1912"
$ws1_data[39,0] = "06:16:15"; $ws1_data[39,1] = "07:58"; $ws1_data[39,2] = "23_HERNANDEZ"; $ws1_data[39,3] = 102; $ws1_data[39,4] = "LP1912"
$ws1_data[40,0] = "06:33:46"; $ws1_data[40,1] = "07:59"; $ws1_data[40,2] = "11_ETCHEVERRY"; $ws1_data[40,3] = 86; $ws1_data[40,4] = "LP1912"
$ws1_data[41,0] = "07:12:53"; $ws1_data[41,1] = "07:59"; $ws1_data[41,2] = "23_HERNANDEZ"; $ws1_data[41,3] = 47; $ws1_data[41,4] = "LP1912"
$ws1_data[42,0] = "06:33:46"; $ws1_data[42,1] = "08:00"; $ws1_data[42,2] = "23_HERNANDEZ"; $ws1_data[42,3] = 87; $ws1_data[42,4] = "LP1912"
$ws1_data[43,0] = "06:16:15"; $ws1_data[43,1] = "08:00"; $ws1_data[43,2] = "11_ETCHEVERRY"; $ws1_data[43,3] = 104; $ws1_data[43,4] = "LP1912"
$ws1_data[44,0] = "07:12:53"; $ws1_data[44,1] = "08:01"; $ws1_data[44,2] = "16_SANTA ANA"; $ws1_data[44,3] = 49; $ws1_data[44,4] = "LP1912"
$ws1_data[45,0] = "06:45:50"; $ws1_data[45,1] = "08:01"; $ws1_data[45,2] = "23_HERNANDEZ"; $ws1_data[45,3] = 76; $ws1_data[45,4] = "LP1912"
$ws1_data[46,0] = "07:36:59"; $ws1_data[46,1] = "08:02"; $ws1_data[46,2] = "17_ROMERO"; $ws1_data[46,3] = 26; $ws1_data[46,4] = "LP1912"
$ws1_data[47,0] = "06:52:52"; $ws1_data[47,1] = "08:02"; $ws1_data[47,2] = "23_HERNANDEZ"; $ws1_data[47,3] = 70; $ws1_data[47,4] = "LP1912"
$ws1_data[48,0] = "06:16:15"; $ws1_data[48,1] = "08:03"; $ws1_data[48,2] = "17_ROMERO"; $ws1_data[48,3] = 107; $ws1_data[48,4] = "LP1912"
$ws1_data[49,0] = "06:33:46"; $ws1_data[49,1] = "08:12"; $ws1_data[49,2] = "10_OLMOS"; $ws1_data[49,3] = 99; $ws1_data[49,4] = "LP1912"
$ws1_data[50,0] = "07:12:53"; $ws1_data[50,1] = "08:13"; $ws1_data[50,2] = "10_OLMOS"; $ws1_data[50,3] = 61; $ws1_data[50,4] = "LP1912"
$ws1_data[51,0] = "07:48:35"; $ws1_data[51,1] = "08:14"; $ws1_data[51,2] = "10_OLMOS"; $ws1_data[51,3] = 26; $ws1_data[51,4] = "LP1912"
$ws1_data[52,0] = "07:36:59"; $ws1_data[52,1] = "08:14"; $ws1_data[52,2] = "17_ROMERO"; $ws1_data[52,3] = 38; $ws1_data[52,4] = "LP1912"
$ws1_data[53,0] = "08:11:27"; $ws1_data[53,1] = "08:14"; $ws1_data[53,2] = "11_ETCHEVERRY"; $ws1_data[53,3] = 3; $ws1_data[53,4] = "LP1912"
$ws1_data[54,0] = "06:16:15"; $ws1_data[54,1] = "08:15"; $ws1_data[54,2] = "17_ROMERO"; $ws1_data[54,3] = 119; $ws1_data[54,4] = "LP1912"
$ws1_data[55,0] = "07:36:59"; $ws1_data[55,1] = "08:25"; $ws1_data[55,2] = "15X38_ABASTO"; $ws1_data[55,3] = 49; $ws1_data[55,4] = "LP1912"
$ws1_data[56,0] = "06:33:46"; $ws1_data[56,1] = "08:26"; $ws1_data[56,2] = "15X38_ABASTO"; $ws1_data[56,3] = 113; $ws1_data[56,4] = "LP1912"
$ws1_data[57,0] = "06:33:46"; $ws1_data[57,1] = "08:27"; $ws1_data[57,2] = "84_COLONIA URQUIZA-ESC 49"; $ws1_data[57,3] = 114; $ws1_data[57,4] = "LP1912"
$ws1_data[58,0] = "06:45:50"; $ws1_data[58,1] = "08:29"; $ws1_data[58,2] = "14_ABASTO"; $ws1_data[58,3] = 104; $ws1_data[58,4] = "LP1912"
$ws1_data[59,0] = "08:29:19"; $ws1_data[59,1] = "08:29"; $ws1_data[59,2] = "23_HERNANDEZ"; $ws1_data[59,3] = 0; $ws1_data[59,4] = "LP1912"
$ws1_data[60,0] = "07:36:59"; $ws1_data[60,1] = "08:30"; $ws1_data[60,2] = "16_P MOR-SANTA ANA"; $ws1_data[60,3] = 54; $ws1_data[60,4] = "LP1912"
$ws1_data[61,0] = "06:33:46"; $ws1_data[61,1] = "08:31"; $ws1_data[61,2] = "16_P MOR-SANTA ANA"; $ws1_data[61,3] = 118; $ws1_data[61,4] = "LP1912"
$ws1_data[62,0] = "08:11:27"; $ws1_data[62,1] = "08:33"; $ws1_data[62,2] = "23_HERNANDEZ"; $ws1_data[62,3] = 22; $ws1_data[62,4] = "LP1912"
$ws1_data[63,0] = "06:45:50"; $ws1_data[63,1] = "08:38"; $ws1_data[63,2] = "215C_EL PATO"; $ws1_data[63,3] = 113; $ws1_data[63,4] = "LP1912"
$ws1_data[64,0] = "07:48:35"; $ws1_data[64,1] = "08:39"; $ws1_data[64,2] = "215C_EL PATO"; $ws1_data[64,3] = 51; $ws1_data[64,4] = "LP1912"
$ws1_data[65,0] = "07:12:53"; $ws1_data[65,1] = "08:43"; $ws1_data[65,2] = "10_OLMOS"; $ws1_data[65,3] = 91; $ws1_data[65,4] = "LP1912"
$ws1_data[66,0] = "07:48:35"; $ws1_data[66,1] = "08:44"; $ws1_data[66,2] = "10_OLMOS"; $ws1_data[66,3] = 56; $ws1_data[66,4] = "LP1912"
$ws1_data[67,0] = "07:12:53"; $ws1_data[67,1] = "08:49"; $ws1_data[67,2] = "215A_EL PATO"; $ws1_data[67,3] = 97; $ws1_data[67,4] = "LP1912"
$ws1_data[68,0] = "07:48:35"; $ws1_data[68,1] = "08:50"; $ws1_data[68,2] = "215A_EL PATO"; $ws1_data[68,3] = 62; $ws1_data[68,4] = "LP1912"
$ws1_data[69,0] = "08:11:27"; $ws1_data[69,1] = "08:53"; $ws1_data[69,2] = "16_SANTA ANA"; $ws1_data[69,3] = 42; $ws1_data[69,4] = "LP1912"
$ws1_data[70,0] = "08:29:19"; $ws1_data[70,1] = "08:54"; $ws1_data[70,2] = "16_SANTA ANA"; $ws1_data[70,3] = 25; $ws1_data[70,4] = "LP1912"
$ws1_data[71,0] = "08:53:12"; $ws1_data[71,1] = "08:55"; $ws1_data[71,2] = "16_SANTA ANA"; $ws1_data[71,3] = 2; $ws1_data[71,4] = "LP1912"
$ws1_data[72,0] = "07:12:53"; $ws1_data[72,1] = "08:59"; $ws1_data[72,2] = "215B_EL PATO"; $ws1_data[72,3] = 107; $ws1_data[72,4] = "LP1912"
$ws1_data[73,0] = "08:39:08"; $ws1_data[73,1] = "09:00"; $ws1_data[73,2] = "215B_EL PATO"; $ws1_data[73,3] = 21; $ws1_data[73,4] = "LP1912"
$ws1_data[74,0] = "07:36:59"; $ws1_data[74,1] = "09:01"; $ws1_data[74,2] = "17X38_ROMERO"; $ws1_data[74,3] = 85; $ws1_data[74,4] = "LP1912"
$ws1_data[75,0] = "07:12:53"; $ws1_data[75,1] = "09:02"; $ws1_data[75,2] = "17X38_ROMERO"; $ws1_data[75,3] = 110; $ws1_data[75,4] = "LP1912"
$ws1_data[76,0] = "07:36:59"; $ws1_data[76,1] = "09:02"; $ws1_data[76,2] = "23_HERNANDEZ"; $ws1_data[76,3] = 86; $ws1_data[76,4] = "LP1912"
$ws1_data[77,0] = "08:29:19"; $ws1_data[77,1] = "09:03"; $ws1_data[77,2] = "23_HERNANDEZ"; $ws1_data[77,3] = 34; $ws1_data[77,4] = "LP1912"
$ws1_data[78,0] = "08:39:08"; $ws1_data[78,1] = "09:04"; $ws1_data[78,2] = "23_HERNANDEZ"; $ws1_data[78,3] = 25; $ws1_data[78,4] = "LP1912"
$ws1_data[79,0] = "07:36:59"; $ws1_data[79,1] = "09:04"; $ws1_data[79,2] = "16_SANTA ANA"; $ws1_data[79,3] = 88; $ws1_data[79,4] = "LP1912"
$ws1_data[80,0] = "07:48:35"; $ws1_data[80,1] = "09:08"; $ws1_data[80,2] = "16_SANTA ANA"; $ws1_data[80,3] = 80; $ws1_data[80,4] = "LP1912"
$ws1_data[81,0] = "07:36:59"; $ws1_data[81,1] = "09:14"; $ws1_data[81,2] = "15_ABASTO"; $ws1_data[81,3] = 98; $ws1_data[81,4] = "LP1912"
$ws1_data[82,0] = "07:36:59"; $ws1_data[82,1] = "09:14"; $ws1_data[82,2] = "11_ETCHEVERRY"; $ws1_data[82,3] = 98; $ws1_data[82,4] = "LP1912"
$ws1_data[83,0] = "08:39:08"; $ws1_data[83,1] = "09:15"; $ws1_data[83,2] = "11_ETCHEVERRY"; $ws1_data[83,3] = 36; $ws1_data[83,4] = "LP1912"
$ws1_data[84,0] = "07:36:59"; $ws1_data[84,1] = "09:16"; $ws1_data[84,2] = "27_EL RETIRO"; $ws1_data[84,3] = 100; $ws1_data[84,4] = "LP1912"
$ws1_data[85,0] = "08:39:08"; $ws1_data[85,1] = "09:17"; $ws1_data[85,2] = "27_EL RETIRO"; $ws1_data[85,3] = 38; $ws1_data[85,4] = "LP1912"
$ws1_data[86,0] = "07:55:46"; $ws1_data[86,1] = "09:21"; $ws1_data[86,2] = "16_SANTA ANA"; $ws1_data[86,3] = 86; $ws1_data[86,4] = "LP1912"
$ws1_data[87,0] = "09:21:49"; $ws1_data[87,1] = "09:22"; $ws1_data[87,2] = "16_SANTA ANA"; $ws1_data[87,3] = 1; $ws1_data[87,4] = "LP1912"
$ws1_data[88,0] = "07:36:59"; $ws1_data[88,1] = "09:26"; $ws1_data[88,2] = "215_EL PELIGRO"; $ws1_data[88,3] = 110; $ws1_data[88,4] = "LP1912"
$ws1_data[89,0] = "07:48:35"; $ws1_data[89,1] = "09:27"; $ws1_data[89,2] = "215_EL PELIGRO"; $ws1_data[89,3] = 99; $ws1_data[89,4] = "LP1912"
$ws1_data[90,0] = "07:36:59"; $ws1_data[90,1] = "09:30"; $ws1_data[90,2] = "16_P MOR-SANTA ANA"; $ws1_data[90,3] = 114; $ws1_data[90,4] = "LP1912"
$ws1_data[91,0] = "07:48:35"; $ws1_data[91,1] = "09:31"; $ws1_data[91,2] = "16_P MOR-SANTA ANA"; $ws1_data[91,3] = 103; $ws1_data[91,4] = "LP1912"
$ws1_data[92,0] = "08:39:08"; $ws1_data[92,1] = "09:32"; $ws1_data[92,2] = "23_HERNANDEZ"; $ws1_data[92,3] = 53; $ws1_data[92,4] = "LP1912"
$ws1_data[93,0] = "08:46:25"; $ws1_data[93,1] = "09:33"; $ws1_data[93,2] = "23_HERNANDEZ"; $ws1_data[93,3] = 47; $ws1_data[93,4] = "LP1912"
$ws1_data[94,0] = "07:48:35"; $ws1_data[94,1] = "09:39"; $ws1_data[94,2] = "15_ABASTO"; $ws1_data[94,3] = 111; $ws1_data[94,4] = "LP1912"
$ws1_data[95,0] = "07:48:35"; $ws1_data[95,1] = "09:44"; $ws1_data[95,2] = "14_ABASTO"; $ws1_data[95,3] = 116; $ws1_data[95,4] = "LP1912"
$ws1_data[96,0] = "08:39:08"; $ws1_data[96,1] = "09:45"; $ws1_data[96,2] = "14_ABASTO"; $ws1_data[96,3] = 66; $ws1_data[96,4] = "LP1912"
$ws1_data[97,0] = "07:55:46"; $ws1_data[97,1] = "09:51"; $ws1_data[97,2] = "15_ABASTO"; $ws1_data[97,3] = 116; $ws1_data[97,4] = "LP1912"
$ws1_data[98,0] = "08:53:12"; $ws1_data[98,1] = "09:54"; $ws1_data[98,2] = "10_OLMOS"; $ws1_data[98,3] = 61; $ws1_data[98,4] = "LP1912"
$ws1_data[99,0] = "09:21:49"; $ws1_data[99,1] = "10:01"; $ws1_data[99,2] = "16_SANTA ANA"; $ws1_data[99,3] = 40; $ws1_data[99,4] = "LP1912"
$ws1_data[100,0] = "09:21:49"; $ws1_data[100,1] = "10:02"; $ws1_data[100,2] = "215C_EL PATO"; $ws1_data[100,3] = 41; $ws1_data[100,4] = "LP1912"
$ws1_data[101,0] = "09:21:49"; $ws1_data[101,1] = "10:03"; $ws1_data[101,2] = "23_HERNANDEZ"; $ws1_data[101,3] = 42; $ws1_data[101,4] = "LP1912"
$ws1_data[102,0] = "08:11:27"; $ws1_data[102,1] = "10:03"; $ws1_data[102,2] = "215C_EL PATO"; $ws1_data[102,3] = 112; $ws1_data[102,4] = "LP1912"
$ws1_data[103,0] = "10:04:17"; $ws1_data[103,1] = "10:04"; $ws1_data[103,2] = "215C_EL PATO"; $ws1_data[103,3] = 0; $ws1_data[103,4] = "LP1912"
$ws1_data[104,0] = "08:46:25"; $ws1_data[104,1] = "10:04"; $ws1_data[104,2] = "14_ABASTO"; $ws1_data[104,3] = 78; $ws1_data[104,4] = "LP1912"
$ws1_data[105,0] = "08:39:08"; $ws1_data[105,1] = "10:05"; $ws1_data[105,2] = "14_ABASTO"; $ws1_data[105,3] = 86; $ws1_data[105,4] = "LP1912"
$ws1_data[106,0] = "10:04:17"; $ws1_data[106,1] = "10:05"; $ws1_data[106,2] = "16_SANTA ANA"; $ws1_data[106,3] = 1; $ws1_data[106,4] = "LP1912"
$ws1_data[107,0] = "08:11:27"; $ws1_data[107,1] = "10:10"; $ws1_data[107,2] = "10_OLMOS"; $ws1_data[107,3] = 119; $ws1_data[107,4] = "LP1912"
$ws1_data[108,0] = "08:29:19"; $ws1_data[108,1] = "10:11"; $ws1_data[108,2] = "27_EL RETIRO"; $ws1_data[108,3] = 102; $ws1_data[108,4] = "LP1912"
$ws1_data[109,0] = "08:29:19"; $ws1_data[109,1] = "10:12"; $ws1_data[109,2] = "10_OLMOS"; $ws1_data[109,3] = 103; $ws1_data[109,4] = "LP1912"
$ws1_data[110,0] = "08:29:19"; $ws1_data[110,1] = "10:14"; $ws1_data[110,2] = "10_OLMOS"; $ws1_data[110,3] = 105; $ws1_data[110,4] = "LP1912"
$ws1_data[111,0] = "09:21:49"; $ws1_data[111,1] = "10:14"; $ws1_data[111,2] = "17_ROMERO"; $ws1_data[111,3] = 53; $ws1_data[111,4] = "LP1912"
$ws1_data[112,0] = "08:39:08"; $ws1_data[112,1] = "10:15"; $ws1_data[112,2] = "10_OLMOS"; $ws1_data[112,3] = 96; $ws1_data[112,4] = "LP1912"
$ws1_data[113,0] = "08:29:19"; $ws1_data[113,1] = "10:15"; $ws1_data[113,2] = "17_ROMERO"; $ws1_data[113,3] = 106; $ws1_data[113,4] = "LP1912"
$ws1_data[114,0] = "08:46:25"; $ws1_data[114,1] = "10:16"; $ws1_data[114,2] = "10_OLMOS"; $ws1_data[114,3] = 90; $ws1_data[114,4] = "LP1912"
$ws1_data[115,0] = "08:46:25"; $ws1_data[115,1] = "10:18"; $ws1_data[115,2] = "27_EL RETIRO"; $ws1_data[115,3] = 92; $ws1_data[115,4] = "LP1912"
$ws1_data[116,0] = "09:21:49"; $ws1_data[116,1] = "10:24"; $ws1_data[116,2] = "11_ETCHEVERRY"; $ws1_data[116,3] = 63; $ws1_data[116,4] = "LP1912"
$ws1_data[117,0] = "10:04:17"; $ws1_data[117,1] = "10:25"; $ws1_data[117,2] = "11_ETCHEVERRY"; $ws1_data[117,3] = 21; $ws1_data[117,4] = "LP1912"
$ws1_data[118,0] = "10:04:17"; $ws1_data[118,1] = "10:25"; $ws1_data[118,2] = "16_SANTA ANA"; $ws1_data[118,3] = 21; $ws1_data[118,4] = "LP1912"
$ws1_data[119,0] = "08:29:19"; $ws1_data[119,1] = "10:26"; $ws1_data[119,2] = "15X38_ABASTO"; $ws1_data[119,3] = 117; $ws1_data[119,4] = "LP1912"
$ws1_data[120,0] = "08:39:08"; $ws1_data[120,1] = "10:30"; $ws1_data[120,2] = "11_ETCHEVERRY"; $ws1_data[120,3] = 111; $ws1_data[120,4] = "LP1912"
$ws1_data[121,0] = "10:04:17"; $ws1_data[121,1] = "10:33"; $ws1_data[121,2] = "23_HERNANDEZ"; $ws1_data[121,3] = 29; $ws1_data[121,4] = "LP1912"
$ws1_data[122,0] = "08:39:08"; $ws1_data[122,1] = "10:34"; $ws1_data[122,2] = "10_OLMOS"; $ws1_data[122,3] = 115; $ws1_data[122,4] = "LP1912"
$ws1_data[123,0] = "10:36:18"; $ws1_data[123,1] = "10:36"; $ws1_data[123,2] = "10_OLMOS"; $ws1_data[123,3] = 0; $ws1_data[123,4] = "LP1912"
$ws1_data[124,0] = "10:36:18"; $ws1_data[124,1] = "10:37"; $ws1_data[124,2] = "16_SANTA ANA"; $ws1_data[124,3] = 1; $ws1_data[124,4] = "LP1912"
$ws1_data[125,0] = "08:39:08"; $ws1_data[125,1] = "10:37"; $ws1_data[125,2] = "16_P MOR-SANTA ANA"; $ws1_data[125,3] = 118; $ws1_data[125,4] = "LP1912"
$ws1_data[126,0] = "10:04:17"; $ws1_data[126,1] = "10:39"; $ws1_data[126,2] = "15_ABASTO"; $ws1_data[126,3] = 35; $ws1_data[126,4] = "LP1912"
$ws1_data[127,0] = "10:04:17"; $ws1_data[127,1] = "10:42"; $ws1_data[127,2] = "27_EL RETIRO"; $ws1_data[127,3] = 38; $ws1_data[127,4] = "LP1912"
$ws1_data[128,0] = "09:21:49"; $ws1_data[128,1] = "10:44"; $ws1_data[128,2] = "14_ABASTO"; $ws1_data[128,3] = 83; $ws1_data[128,4] = "LP1912"
$ws1_data[129,0] = "10:04:17"; $ws1_data[129,1] = "10:45"; $ws1_data[129,2] = "14_ABASTO"; $ws1_data[129,3] = 41; $ws1_data[129,4] = "LP1912"
$ws1_data[130,0] = "10:48:14"; $ws1_data[130,1] = "10:48"; $ws1_data[130,2] = "10_OLMOS"; $ws1_data[130,3] = 0; $ws1_data[130,4] = "LP1912"
$ws1_data[131,0] = "10:48:14"; $ws1_data[131,1] = "10:49"; $ws1_data[131,2] = "16_SANTA ANA"; $ws1_data[131,3] = 1; $ws1_data[131,4] = "LP1912"
$ws1_data[132,0] = "10:04:17"; $ws1_data[132,1] = "10:51"; $ws1_data[132,2] = "15_ABASTO"; $ws1_data[132,3] = 47; $ws1_data[132,4] = "LP1912"
$ws1_data[133,0] = "10:36:18"; $ws1_data[133,1] = "10:54"; $ws1_data[133,2] = "10_OLMOS"; $ws1_data[133,3] = 18; $ws1_data[133,4] = "LP1912"
$ws1_data[134,0] = "10:55:25"; $ws1_data[134,1] = "10:56"; $ws1_data[134,2] = "10_OLMOS"; $ws1_data[134,3] = 1; $ws1_data[134,4] = "LP1912"
$ws1_data[135,0] = "10:55:25"; $ws1_data[135,1] = "10:56"; $ws1_data[135,2] = "16_SANTA ANA"; $ws1_data[135,3] = 1; $ws1_data[135,4] = "LP1912"
$ws1_data[136,0] = "09:21:49"; $ws1_data[136,1] = "10:56"; $ws1_data[136,2] = "27_EL RETIRO"; $ws1_data[136,3] = 95; $ws1_data[136,4] = "LP1912"
$ws1_data[137,0] = "10:04:17"; $ws1_data[137,1] = "10:57"; $ws1_data[137,2] = "27_EL RETIRO"; $ws1_data[137,3] = 53; $ws1_data[137,4] = "LP1912"
$ws1_data[138,0] = "10:55:25"; $ws1_data[138,1] = "10:59"; $ws1_data[138,2] = "27_EL RETIRO"; $ws1_data[138,3] = 4; $ws1_data[138,4] = "LP1912"
$ws1_data[139,0] = "09:21:49"; $ws1_data[139,1] = "11:01"; $ws1_data[139,2] = "17_ROMERO"; $ws1_data[139,3] = 100; $ws1_data[139,4] = "LP1912"
$ws1_data[140,0] = "10:36:18"; $ws1_data[140,1] = "11:03"; $ws1_data[140,2] = "23_HERNANDEZ"; $ws1_data[140,3] = 27; $ws1_data[140,4] = "LP1912"
$ws1_data[141,0] = "09:21:49"; $ws1_data[141,1] = "11:04"; $ws1_data[141,2] = "14_ABASTO"; $ws1_data[141,3] = 103; $ws1_data[141,4] = "LP1912"
$ws1_data[142,0] = "10:04:17"; $ws1_data[142,1] = "11:05"; $ws1_data[142,2] = "14_ABASTO"; $ws1_data[142,3] = 61; $ws1_data[142,4] = "LP1912"
$ws1_data[143,0] = "10:36:18"; $ws1_data[143,1] = "11:11"; $ws1_data[143,2] = "15_ABASTO"; $ws1_data[143,3] = 35; $ws1_data[143,4] = "LP1912"
$ws1_data[144,0] = "11:11:31"; $ws1_data[144,1] = "11:11"; $ws1_data[144,2] = "16_SANTA ANA"; $ws1_data[144,3] = 0; $ws1_data[144,4] = "LP1912"
$ws1_data[145,0] = "10:04:17"; $ws1_data[145,1] = "11:11"; $ws1_data[145,2] = "23_HERNANDEZ"; $ws1_data[145,3] = 67; $ws1_data[145,4] = "LP1912"
$ws1_data[146,0] = "09:21:49"; $ws1_data[146,1] = "11:14"; $ws1_data[146,2] = "225_C ROCA-H SUR"; $ws1_data[146,3] = 113; $ws1_data[146,4] = "LP1912"
$ws1_data[147,0] = "09:21:49"; $ws1_data[147,1] = "11:20"; $ws1_data[147,2] = "215C_EL PATO"; $ws1_data[147,3] = 119; $ws1_data[147,4] = "LP1912"
$ws1_data[148,0] = "10:48:14"; $ws1_data[148,1] = "11:21"; $ws1_data[148,2] = "10_OLMOS"; $ws1_data[148,3] = 33; $ws1_data[148,4] = "LP1912"
$ws1_data[149,0] = "10:04:17"; $ws1_data[149,1] = "11:21"; $ws1_data[149,2] = "215C_EL PATO"; $ws1_data[149,3] = 77; $ws1_data[149,4] = "LP1912"
$ws1_data[150,0] = "10:36:18"; $ws1_data[150,1] = "11:22"; $ws1_data[150,2] = "10_OLMOS"; $ws1_data[150,3] = 46; $ws1_data[150,4] = "LP1912"
$ws1_data[151,0] = "10:36:18"; $ws1_data[151,1] = "11:24"; $ws1_data[151,2] = "11_ETCHEVERRY"; $ws1_data[151,3] = 48; $ws1_data[151,4] = "LP1912"
$ws1_data[152,0] = "10:36:18"; $ws1_data[152,1] = "11:25"; $ws1_data[152,2] = "16_P MOR-SANTA ANA"; $ws1_data[152,3] = 49; $ws1_data[152,4] = "LP1912"
$ws1_data[153,0] = "10:04:17"; $ws1_data[153,1] = "11:25"; $ws1_data[153,2] = "11_ETCHEVERRY"; $ws1_data[153,3] = 81; $ws1_data[153,4] = "LP1912"
$ws1_data[154,0] = "10:04:17"; $ws1_data[154,1] = "11:30"; $ws1_data[154,2] = "15X38_ABASTO"; $ws1_data[154,3] = 86; $ws1_data[154,4] = "LP1912"
$ws1_data[155,0] = "10:48:14"; $ws1_data[155,1] = "11:32"; $ws1_data[155,2] = "23_HERNANDEZ"; $ws1_data[155,3] = 44; $ws1_data[155,4] = "LP1912"
$ws1_data[156,0] = "10:36:18"; $ws1_data[156,1] = "11:33"; $ws1_data[156,2] = "23_HERNANDEZ"; $ws1_data[156,3] = 57; $ws1_data[156,4] = "LP1912"
$ws1_data[157,0] = "10:48:14"; $ws1_data[157,1] = "11:33"; $ws1_data[157,2] = "16_SANTA ANA"; $ws1_data[157,3] = 45; $ws1_data[157,4] = "LP1912"
$ws1_data[158,0] = "10:04:17"; $ws1_data[158,1] = "11:34"; $ws1_data[158,2] = "10_OLMOS"; $ws1_data[158,3] = 90; $ws1_data[158,4] = "LP1912"
$ws1_data[159,0] = "11:34:25"; $ws1_data[159,1] = "11:34"; $ws1_data[159,2] = "23_HERNANDEZ"; $ws1_data[159,3] = 0; $ws1_data[159,4] = "LP1912"
$ws1_data[160,0] = "11:34:25"; $ws1_data[160,1] = "11:34"; $ws1_data[160,2] = "16_SANTA ANA"; $ws1_data[160,3] = 0; $ws1_data[160,4] = "LP1912"
$ws1_data[161,0] = "10:36:18"; $ws1_data[161,1] = "11:35"; $ws1_data[161,2] = "16_SANTA ANA"; $ws1_data[161,3] = 59; $ws1_data[161,4] = "LP1912"
$ws1_data[162,0] = "10:04:17"; $ws1_data[162,1] = "11:37"; $ws1_data[162,2] = "16_P MOR-SANTA ANA"; $ws1_data[162,3] = 93; $ws1_data[162,4] = "LP1912"
$ws1_data[163,0] = "10:04:17"; $ws1_data[163,1] = "11:40"; $ws1_data[163,2] = "215A_EL PATO"; $ws1_data[163,3] = 96; $ws1_data[163,4] = "LP1912"
$ws1_data[164,0] = "11:34:25"; $ws1_data[164,1] = "11:41"; $ws1_data[164,2] = "10_OLMOS"; $ws1_data[164,3] = 7; $ws1_data[164,4] = "LP1912"
$ws1_data[165,0] = "10:55:25"; $ws1_data[165,1] = "11:44"; $ws1_data[165,2] = "215B_EL PATO"; $ws1_data[165,3] = 49; $ws1_data[165,4] = "LP1912"
$ws1_data[166,0] = "10:04:17"; $ws1_data[166,1] = "11:45"; $ws1_data[166,2] = "215B_EL PATO"; $ws1_data[166,3] = 101; $ws1_data[166,4] = "LP1912"
$ws1_data[167,0] = "10:55:25"; $ws1_data[167,1] = "11:53"; $ws1_data[167,2] = "15_ABASTO"; $ws1_data[167,3] = 58; $ws1_data[167,4] = "LP1912"
$ws1_data[168,0] = "11:53:59"; $ws1_data[168,1] = "11:53"; $ws1_data[168,2] = "225_GOMEZ"; $ws1_data[168,3] = 0; $ws1_data[168,4] = "LP1912"
$ws1_data[169,0] = "10:04:17"; $ws1_data[169,1] = "11:54"; $ws1_data[169,2] = "225_GOMEZ"; $ws1_data[169,3] = 110; $ws1_data[169,4] = "LP1912"
$ws1_data[170,0] = "11:47:13"; $ws1_data[170,1] = "11:57"; $ws1_data[170,2] = "16_SANTA ANA"; $ws1_data[170,3] = 10; $ws1_data[170,4] = "LP1912"
$ws1_data[171,0] = "11:34:25"; $ws1_data[171,1] = "12:03"; $ws1_data[171,2] = "23_HERNANDEZ"; $ws1_data[171,3] = 29; $ws1_data[171,4] = "LP1912"
$ws1_data[172,0] = "11:53:59"; $ws1_data[172,1] = "12:04"; $ws1_data[172,2] = "17_ROMERO"; $ws1_data[172,3] = 11; $ws1_data[172,4] = "LP1912"
$ws1_data[173,0] = "11:11:31"; $ws1_data[173,1] = "12:05"; $ws1_data[173,2] = "17_ROMERO"; $ws1_data[173,3] = 54; $ws1_data[173,4] = "LP1912"
$ws1_data[174,0] = "11:11:31"; $ws1_data[174,1] = "12:06"; $ws1_data[174,2] = "23_HERNANDEZ"; $ws1_data[174,3] = 55; $ws1_data[174,4] = "LP1912"
$ws1_data[175,0] = "10:48:14"; $ws1_data[175,1] = "12:07"; $ws1_data[175,2] = "14_ABASTO"; $ws1_data[175,3] = 79; $ws1_data[175,4] = "LP1912"
$ws1_data[176,0] = "11:47:13"; $ws1_data[176,1] = "12:07"; $ws1_data[176,2] = "23_HERNANDEZ"; $ws1_data[176,3] = 20; $ws1_data[176,4] = "LP1912"
$ws1_data[177,0] = "12:11:45"; $ws1_data[177,1] = "12:11"; $ws1_data[177,2] = "16_SANTA ANA"; $ws1_data[177,3] = 0; $ws1_data[177,4] = "LP1912"
$ws1_data[178,0] = "12:11:45"; $ws1_data[178,1] = "12:11"; $ws1_data[178,2] = "15_ABASTO"; $ws1_data[178,3] = 0; $ws1_data[178,4] = "LP1912"
$ws1_data[179,0] = "11:53:59"; $ws1_data[179,1] = "12:16"; $ws1_data[179,2] = "15_ABASTO"; $ws1_data[179,3] = 23; $ws1_data[179,4] = "LP1912"
$ws1_data[180,0] = "11:53:59"; $ws1_data[180,1] = "12:17"; $ws1_data[180,2] = "10_OLMOS"; $ws1_data[180,3] = 24; $ws1_data[180,4] = "LP1912"
$ws1_data[181,0] = "12:11:45"; $ws1_data[181,1] = "12:17"; $ws1_data[181,2] = "27_EL RETIRO"; $ws1_data[181,3] = 6; $ws1_data[181,4] = "LP1912"
$ws1_data[182,0] = "11:11:31"; $ws1_data[182,1] = "12:17"; $ws1_data[182,2] = "15_ABASTO"; $ws1_data[182,3] = 66; $ws1_data[182,4] = "LP1912"
$ws1_data[183,0] = "11:11:31"; $ws1_data[183,1] = "12:18"; $ws1_data[183,2] = "10_OLMOS"; $ws1_data[183,3] = 67; $ws1_data[183,4] = "LP1912"
$ws1_data[184,0] = "11:34:25"; $ws1_data[184,1] = "12:20"; $ws1_data[184,2] = "17_ROMERO"; $ws1_data[184,3] = 46; $ws1_data[184,4] = "LP1912"
$ws1_data[185,0] = "11:47:13"; $ws1_data[185,1] = "12:21"; $ws1_data[185,2] = "17_ROMERO"; $ws1_data[185,3] = 34; $ws1_data[185,4] = "LP1912"
$ws1_data[186,0] = "11:53:59"; $ws1_data[186,1] = "12:28"; $ws1_data[186,2] = "215C_EL PATO"; $ws1_data[186,3] = 35; $ws1_data[186,4] = "LP1912"
$ws1_data[187,0] = "10:36:18"; $ws1_data[187,1] = "12:29"; $ws1_data[187,2] = "215C_EL PATO"; $ws1_data[187,3] = 113; $ws1_data[187,4] = "LP1912"
$ws1_data[188,0] = "10:36:18"; $ws1_data[188,1] = "12:30"; $ws1_data[188,2] = "11_ETCHEVERRY"; $ws1_data[188,3] = 114; $ws1_data[188,4] = "LP1912"
$ws1_data[189,0] = "11:53:59"; $ws1_data[189,1] = "12:30"; $ws1_data[189,2] = "16_P MOR-SANTA ANA"; $ws1_data[189,3] = 37; $ws1_data[189,4] = "LP1912"
$ws1_data[190,0] = "10:48:14"; $ws1_data[190,1] = "12:31"; $ws1_data[190,2] = "11_ETCHEVERRY"; $ws1_data[190,3] = 103; $ws1_data[190,4] = "LP1912"
$ws1_data[191,0] = "10:36:18"; $ws1_data[191,1] = "12:31"; $ws1_data[191,2] = "16_P MOR-SANTA ANA"; $ws1_data[191,3] = 115; $ws1_data[191,4] = "LP1912"
$ws1_data[192,0] = "12:32:47"; $ws1_data[192,1] = "12:32"; $ws1_data[192,2] = "16_P MOR-SANTA ANA"; $ws1_data[192,3] = 0; $ws1_data[192,4] = "LP1912"
$ws1_data[193,0] = "12:32:47"; $ws1_data[193,1] = "12:33"; $ws1_data[193,2] = "27_EL RETIRO"; $ws1_data[193,3] = 1; $ws1_data[193,4] = "LP1912"
$ws1_data[194,0] = "12:32:47"; $ws1_data[194,1] = "12:33"; $ws1_data[194,2] = "17_ROMERO"; $ws1_data[194,3] = 1; $ws1_data[194,4] = "LP1912"
$ws1_data[195,0] = "10:55:25"; $ws1_data[195,1] = "12:36"; $ws1_data[195,2] = "27_EL RETIRO"; $ws1_data[195,3] = 101; $ws1_data[195,4] = "LP1912"
$ws1_data[196,0] = "10:48:14"; $ws1_data[196,1] = "12:37"; $ws1_data[196,2] = "27_EL RETIRO"; $ws1_data[196,3] = 109; $ws1_data[196,4] = "LP1912"
$ws1_data[197,0] = "11:53:59"; $ws1_data[197,1] = "12:39"; $ws1_data[197,2] = "15X38_ABASTO"; $ws1_data[197,3] = 46; $ws1_data[197,4] = "LP1912"
$ws1_data[198,0] = "10:48:14"; $ws1_data[198,1] = "12:40"; $ws1_data[198,2] = "15X38_ABASTO"; $ws1_data[198,3] = 112; $ws1_data[198,4] = "LP1912"
$ws1_data[199,0] = "10:55:25"; $ws1_data[199,1] = "12:42"; $ws1_data[199,2] = "14_ABASTO"; $ws1_data[199,3] = 107; $ws1_data[199,4] = "LP1912"
$ws1_data[200,0] = "10:55:25"; $ws1_data[200,1] = "12:43"; $ws1_data[200,2] = "15X38_ABASTO"; $ws1_data[200,3] = 108; $ws1_data[200,4] = "LP1912"
$ws1_data[201,0] = "10:48:14"; $ws1_data[201,1] = "12:43"; $ws1_data[201,2] = "14_ABASTO"; $ws1_data[201,3] = 115; $ws1_data[201,4] = "LP1912"
$ws1_data[202,0] = "11:53:59"; $ws1_data[202,1] = "12:50"; $ws1_data[202,2] = "15_ABASTO"; $ws1_data[202,3] = 57; $ws1_data[202,4] = "LP1912"
$ws1_data[203,0] = "11:47:13"; $ws1_data[203,1] = "12:51"; $ws1_data[203,2] = "15_ABASTO"; $ws1_data[203,3] = 64; $ws1_data[203,4] = "LP1912"
$ws1_data[204,0] = "11:11:31"; $ws1_data[204,1] = "12:54"; $ws1_data[204,2] = "15X38_ABASTO"; $ws1_data[204,3] = 103; $ws1_data[204,4] = "LP1912"
$ws1_data[205,0] = "12:32:47"; $ws1_data[205,1] = "12:54"; $ws1_data[205,2] = "10_OLMOS"; $ws1_data[205,3] = 22; $ws1_data[205,4] = "LP1912"
$ws1_data[206,0] = "11:53:59"; $ws1_data[206,1] = "13:00"; $ws1_data[206,2] = "215C_EL PATO"; $ws1_data[206,3] = 67; $ws1_data[206,4] = "LP1912"
$ws1_data[207,0] = "11:11:31"; $ws1_data[207,1] = "13:01"; $ws1_data[207,2] = "215C_EL PATO"; $ws1_data[207,3] = 110; $ws1_data[207,4] = "LP1912"
$ws1_data[208,0] = "12:32:47"; $ws1_data[208,1] = "13:03"; $ws1_data[208,2] = "23_HERNANDEZ"; $ws1_data[208,3] = 31; $ws1_data[208,4] = "LP1912"
$ws1_data[209,0] = "11:47:13"; $ws1_data[209,1] = "13:05"; $ws1_data[209,2] = "23_HERNANDEZ"; $ws1_data[209,3] = 78; $ws1_data[209,4] = "LP1912"
$ws1_data[210,0] = "12:11:45"; $ws1_data[210,1] = "13:06"; $ws1_data[210,2] = "23_HERNANDEZ"; $ws1_data[210,3] = 55; $ws1_data[210,4] = "LP1912"
$ws1_data[211,0] = "11:11:31"; $ws1_data[211,1] = "13:06"; $ws1_data[211,2] = "14_ABASTO"; $ws1_data[211,3] = 115; $ws1_data[211,4] = "LP1912"
$ws1_data[212,0] = "11:47:13"; $ws1_data[212,1] = "13:07"; $ws1_data[212,2] = "14_ABASTO"; $ws1_data[212,3] = 80; $ws1_data[212,4] = "LP1912"
$ws1_data[213,0] = "11:53:59"; $ws1_data[213,1] = "13:07"; $ws1_data[213,2] = "23_HERNANDEZ"; $ws1_data[213,3] = 74; $ws1_data[213,4] = "LP1912"
$ws1_data[214,0] = "11:53:59"; $ws1_data[214,1] = "13:10"; $ws1_data[214,2] = "16_SANTA ANA"; $ws1_data[214,3] = 77; $ws1_data[214,4] = "LP1912"
$ws1_data[215,0] = "11:53:59"; $ws1_data[215,1] = "13:10"; $ws1_data[215,2] = "215_ALUAR"; $ws1_data[215,3] = 77; $ws1_data[215,4] = "LP1912"
$ws1_data[216,0] = "11:34:25"; $ws1_data[216,1] = "13:11"; $ws1_data[216,2] = "215_ALUAR"; $ws1_data[216,3] = 97; $ws1_data[216,4] = "LP1912"
$ws1_data[217,0] = "11:47:13"; $ws1_data[217,1] = "13:11"; $ws1_data[217,2] = "16_SANTA ANA"; $ws1_data[217,3] = 84; $ws1_data[217,4] = "LP1912"
$ws1_data[218,0] = "11:34:25"; $ws1_data[218,1] = "13:18"; $ws1_data[218,2] = "11_ETCHEVERRY"; $ws1_data[218,3] = 104; $ws1_data[218,4] = "LP1912"
$ws1_data[219,0] = "11:47:13"; $ws1_data[219,1] = "13:19"; $ws1_data[219,2] = "11_ETCHEVERRY"; $ws1_data[219,3] = 92; $ws1_data[219,4] = "LP1912"
$ws1_data[220,0] = "11:53:59"; $ws1_data[220,1] = "13:20"; $ws1_data[220,2] = "16_SANTA ANA"; $ws1_data[220,3] = 87; $ws1_data[220,4] = "LP1912"
$ws1_data[221,0] = "11:53:59"; $ws1_data[221,1] = "13:20"; $ws1_data[221,2] = "17_ROMERO"; $ws1_data[221,3] = 87; $ws1_data[221,4] = "LP1912"
$ws1_data[222,0] = "12:32:47"; $ws1_data[222,1] = "13:21"; $ws1_data[222,2] = "16_SANTA ANA"; $ws1_data[222,3] = 49; $ws1_data[222,4] = "LP1912"
$ws1_data[223,0] = "11:34:25"; $ws1_data[223,1] = "13:21"; $ws1_data[223,2] = "17_ROMERO"; $ws1_data[223,3] = 107; $ws1_data[223,4] = "LP1912"
$ws1_data[224,0] = "11:53:59"; $ws1_data[224,1] = "13:29"; $ws1_data[224,2] = "10_OLMOS"; $ws1_data[224,3] = 96; $ws1_data[224,4] = "LP1912"
$ws1_data[225,0] = "11:53:59"; $ws1_data[225,1] = "13:29"; $ws1_data[225,2] = "215A_EL PATO"; $ws1_data[225,3] = 96; $ws1_data[225,4] = "LP1912"
$ws1_data[226,0] = "11:34:25"; $ws1_data[226,1] = "13:30"; $ws1_data[226,2] = "215A_EL PATO"; $ws1_data[226,3] = 116; $ws1_data[226,4] = "LP1912"
$ws1_data[227,0] = "11:47:13"; $ws1_data[227,1] = "13:30"; $ws1_data[227,2] = "10_OLMOS"; $ws1_data[227,3] = 103; $ws1_data[227,4] = "LP1912"
$ws1_data[228,0] = "11:53:59"; $ws1_data[228,1] = "13:30"; $ws1_data[228,2] = "16_P MOR-SANTA ANA"; $ws1_data[228,3] = 97; $ws1_data[228,4] = "LP1912"
$ws1_data[229,0] = "11:47:13"; $ws1_data[229,1] = "13:31"; $ws1_data[229,2] = "16_P MOR-SANTA ANA"; $ws1_data[229,3] = 104; $ws1_data[229,4] = "LP1912"
$ws1_data[230,0] = "12:11:45"; $ws1_data[230,1] = "13:37"; $ws1_data[230,2] = "23_HERNANDEZ"; $ws1_data[230,3] = 86; $ws1_data[230,4] = "LP1912"
$ws1_data[231,0] = "11:53:59"; $ws1_data[231,1] = "13:39"; $ws1_data[231,2] = "17X38_ROMERO"; $ws1_data[231,3] = 106; $ws1_data[231,4] = "LP1912"
$ws1_data[232,0] = "12:32:47"; $ws1_data[232,1] = "13:39"; $ws1_data[232,2] = "23_HERNANDEZ"; $ws1_data[232,3] = 67; $ws1_data[232,4] = "LP1912"
$ws1_data[233,0] = "11:47:13"; $ws1_data[233,1] = "13:40"; $ws1_data[233,2] = "17X38_ROMERO"; $ws1_data[233,3] = 113; $ws1_data[233,4] = "LP1912"
$ws1_data[234,0] = "12:32:47"; $ws1_data[234,1] = "13:51"; $ws1_data[234,2] = "15_ABASTO"; $ws1_data[234,3] = 79; $ws1_data[234,4] = "LP1912"
$ws1_data[235,0] = "12:11:45"; $ws1_data[235,1] = "13:54"; $ws1_data[235,2] = "225_GOMEZ"; $ws1_data[235,3] = 103; $ws1_data[235,4] = "LP1912"
$ws1_data[236,0] = "12:11:45"; $ws1_data[236,1] = "14:01"; $ws1_data[236,2] = "17_ROMERO"; $ws1_data[236,3] = 110; $ws1_data[236,4] = "LP1912"
$ws1_data[237,0] = "12:32:47"; $ws1_data[237,1] = "14:20"; $ws1_data[237,2] = "215C_EL PATO"; $ws1_data[237,3] = 108; $ws1_data[237,4] = "LP1912"
$ws1_data[238,0] = "12:32:47"; $ws1_data[238,1] = "14:24"; $ws1_data[238,2] = "11_ETCHEVERRY"; $ws1_data[238,3] = 112; $ws1_data[238,4] = "LP1912"
$ws1.Range("A6:E244").Value = $ws1_data

# ---- Worksheet: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 12:32:47"
$ws2.Range("A3").Value = "Total filas: 33"
$ws2_data = New-Object "object[,]" 33,5
$ws2_data[0,0] = "05:42:22"; $ws2_data[0,1] = "06:11"; $ws2_data[0,2] = "215A_EL PATO"; $ws2_data[0,3] = 29; $ws2_data[0,4] = "LP1912"
$ws2_data[1,0] = "06:33:46"; $ws2_data[1,1] = "06:50"; $ws2_data[1,2] = "215A_EL PATO"; $ws2_data[1,3] = 17; $ws2_data[1,4] = "LP1912"
$ws2_data[2,0] = "06:16:15"; $ws2_data[2,1] = "06:51"; $ws2_data[2,2] = "215A_EL PATO"; $ws2_data[2,3] = 35; $ws2_data[2,4] = "LP1912"
$ws2_data[3,0] = "06:52:52"; $ws2_data[3,1] = "06:52"; $ws2_data[3,2] = "215A_EL PATO"; $ws2_data[3,3] = 0; $ws2_data[3,4] = "LP1912"
$ws2_data[4,0] = "06:33:46"; $ws2_data[4,1] = "07:06"; $ws2_data[4,2] = "215C_EL PATO"; $ws2_data[4,3] = 33; $ws2_data[4,4] = "LP1912"
$ws2_data[5,0] = "05:42:22"; $ws2_data[5,1] = "07:07"; $ws2_data[5,2] = "215C_EL PATO"; $ws2_data[5,3] = 85; $ws2_data[5,4] = "LP1912"
$ws2_data[6,0] = "06:33:46"; $ws2_data[6,1] = "07:20"; $ws2_data[6,2] = "215A_EL PATO"; $ws2_data[6,3] = 47; $ws2_data[6,4] = "LP1912"
$ws2_data[7,0] = "05:42:22"; $ws2_data[7,1] = "07:21"; $ws2_data[7,2] = "215A_EL PATO"; $ws2_data[7,3] = 99; $ws2_data[7,4] = "LP1912"
$ws2_data[8,0] = "06:45:50"; $ws2_data[8,1] = "08:38"; $ws2_data[8,2] = "215C_EL PATO"; $ws2_data[8,3] = 113; $ws2_data[8,4] = "LP1912"
$ws2_data[9,0] = "07:48:35"; $ws2_data[9,1] = "08:39"; $ws2_data[9,2] = "215C_EL PATO"; $ws2_data[9,3] = 51; $ws2_data[9,4] = "LP1912"
$ws2_data[10,0] = "07:12:53"; $ws2_data[10,1] = "08:49"; $ws2_data[10,2] = "215A_EL PATO"; $ws2_data[10,3] = 97; $ws2_data[10,4] = "LP1912"
$ws2_data[11,0] = "07:48:35"; $ws2_data[11,1] = "08:50"; $ws2_data[11,2] = "215A_EL PATO"; $ws2_data[11,3] = 62; $ws2_data[11,4] = "LP1912"
$ws2_data[12,0] = "07:12:53"; $ws2_data[12,1] = "08:59"; $ws2_data[12,2] = "215B_EL PATO"; $ws2_data[12,3] = 107; $ws2_data[12,4] = "LP1912"
$ws2_data[13,0] = "08:39:08"; $ws2_data[13,1] = "09:00"; $ws2_data[13,2] = "215B_EL PATO"; $ws2_data[13,3] = 21; $ws2_data[13,4] = "LP1912"
$ws2_data[14,0] = "07:36:59"; $ws2_data[14,1] = "09:26"; $ws2_data[14,2] = "215_EL PELIGRO"; $ws2_data[14,3] = 110; $ws2_data[14,4] = "LP1912"
$ws2_data[15,0] = "07:48:35"; $ws2_data[15,1] = "09:27"; $ws2_data[15,2] = "215_EL PELIGRO"; $ws2_data[15,3] = 99; $ws2_data[15,4] = "LP1912"
$ws2_data[16,0] = "09:21:49"; $ws2_data[16,1] = "10:02"; $ws2_data[16,2] = "215C_EL PATO"; $ws2_data[16,3] = 41; $ws2_data[16,4] = "LP1912"
$ws2_data[17,0] = "08:11:27"; $ws2_data[17,1] = "10:03"; $ws2_data[17,2] = "215C_EL PATO"; $ws2_data[17,3] = 112; $ws2_data[17,4] = "LP1912"
$ws2_data[18,0] = "10:04:17"; $ws2_data[18,1] = "10:04"; $ws2_data[18,2] = "215C_EL PATO"; $ws2_data[18,3] = 0; $ws2_data[18,4] = "LP1912"
$ws2_data[19,0] = "09:21:49"; $ws2_data[19,1] = "11:20"; $ws2_data[19,2] = "215C_EL PATO"; $ws2_data[19,3] = 119; $ws2_data[19,4] = "LP1912"
$ws2_data[20,0] = "10:04:17"; $ws2_data[20,1] = "11:21"; $ws2_data[20,2] = "215C_EL PATO"; $ws2_data[20,3] = 77; $ws2_data[20,4] = "LP1912"
$ws2_data[21,0] = "10:04:17"; $ws2_data[21,1] = "11:40"; $ws2_data[21,2] = "215A_EL PATO"; $ws2_data[21,3] = 96; $ws2_data[21,4] = "LP1912"
$ws2_data[22,0] = "10:55:25"; $ws2_data[22,1] = "11:44"; $ws2_data[22,2] = "215B_EL PATO"; $ws2_data[22,3] = 49; $ws2_data[22,4] = "LP1912"
$ws2_data[23,0] = "10:04:17"; $ws2_data[23,1] = "11:45"; $ws2_data[23,2] = "215B_EL PATO"; $ws2_data[23,3] = 101; $ws2_data[23,4] = "LP1912"
$ws2_data[24,0] = "11:53:59"; $ws2_data[24,1] = "12:28"; $ws2_data[24,2] = "215C_EL PATO"; $ws2_data[24,3] = 35; $ws2_data[24,4] = "LP1912"
$ws2_data[25,0] = "10:36:18"; $ws2_data[25,1] = "12:29"; $ws2_data[25,2] = "215C_EL PATO"; $ws2_data[25,3] = 113; $ws2_data[25,4] = "LP1912"
$ws2_data[26,0] = "11:53:59"; $ws2_data[26,1] = "13:00"; $ws2_data[26,2] = "215C_EL PATO"; $ws2_data[26,3] = 67; $ws2_data[26,4] = "LP1912"
$ws2_data[27,0] = "11:11:31"; $ws2_data[27,1] = "13:01"; $ws2_data[27,2] = "215C_EL PATO"; $ws2_data[27,3] = 110; $ws2_data[27,4] = "LP1912"
$ws2_data[28,0] = "11:53:59"; $ws2_data[28,1] = "13:10"; $ws2_data[28,2] = "215_ALUAR"; $ws2_data[28,3] = 77; $ws2_data[28,4] = "LP1912"
$ws2_data[29,0] = "11:34:25"; $ws2_data[29,1] = "13:11"; $ws2_data[29,2] = "215_ALUAR"; $ws2_data[29,3] = 97; $ws2_data[29,4] = "LP1912"
$ws2_data[30,0] = "11:53:59"; $ws2_data[30,1] = "13:29"; $ws2_data[30,2] = "215A_EL PATO"; $ws2_data[30,3] = 96; $ws2_data[30,4] = "LP1912"
$ws2_data[31,0] = "11:34:25"; $ws2_data[31,1] = "13:30"; $ws2_data[31,2] = "215A_EL PATO"; $ws2_data[31,3] = 116; $ws2_data[31,4] = "LP1912"
$ws2_data[32,0] = "12:32:47"; $ws2_data[32,1] = "14:20"; $ws2_data[32,2] = "215C_EL PATO"; $ws2_data[32,3] = 108; $ws2_data[32,4] = "LP1912"
$ws2.Range("A6:E38").Value = $ws2_data

# ---- Worksheet: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 12:32:47"
$ws3.Range("A3").Value = "Total filas: 26"
$ws3_data = New-Object "object[,]" 26,5
$ws3_data[0,0] = "05:42:22"; $ws3_data[0,1] = "07:27"; $ws3_data[0,2] = "215A_LA PLATA"; $ws3_data[0,3] = 105; $ws3_data[0,4] = "L6173"
$ws3_data[1,0] = "06:33:46"; $ws3_data[1,1] = "08:09"; $ws3_data[1,2] = "215A_LA PLATA"; $ws3_data[1,3] = 96; $ws3_data[1,4] = "L6173"
$ws3_data[2,0] = "06:16:15"; $ws3_data[2,1] = "08:10"; $ws3_data[2,2] = "215A_LA PLATA"; $ws3_data[2,3] = 114; $ws3_data[2,4] = "L6173"
$ws3_data[3,0] = "08:11:27"; $ws3_data[3,1] = "08:11"; $ws3_data[3,2] = "215A_LA PLATA"; $ws3_data[3,3] = 0; $ws3_data[3,4] = "L6173"
$ws3_data[4,0] = "06:33:46"; $ws3_data[4,1] = "08:22"; $ws3_data[4,2] = "215C_LA PLATA"; $ws3_data[4,3] = 109; $ws3_data[4,4] = "L6203"
$ws3_data[5,0] = "07:48:35"; $ws3_data[5,1] = "08:25"; $ws3_data[5,2] = "215C_LA PLATA"; $ws3_data[5,3] = 37; $ws3_data[5,4] = "L6203"
$ws3_data[6,0] = "07:55:46"; $ws3_data[6,1] = "08:26"; $ws3_data[6,2] = "215C_LA PLATA"; $ws3_data[6,3] = 31; $ws3_data[6,4] = "L6203"
$ws3_data[7,0] = "07:36:59"; $ws3_data[7,1] = "08:27"; $ws3_data[7,2] = "215C_LA PLATA"; $ws3_data[7,3] = 51; $ws3_data[7,4] = "L6203"
$ws3_data[8,0] = "08:46:25"; $ws3_data[8,1] = "08:48"; $ws3_data[8,2] = "215A_LA PLATA"; $ws3_data[8,3] = 2; $ws3_data[8,4] = "L6173"
$ws3_data[9,0] = "07:36:59"; $ws3_data[9,1] = "08:51"; $ws3_data[9,2] = "215A_LA PLATA"; $ws3_data[9,3] = 75; $ws3_data[9,4] = "L6173"
$ws3_data[10,0] = "07:48:35"; $ws3_data[10,1] = "08:52"; $ws3_data[10,2] = "215A_LA PLATA"; $ws3_data[10,3] = 64; $ws3_data[10,4] = "L6173"
$ws3_data[11,0] = "08:53:12"; $ws3_data[11,1] = "08:53"; $ws3_data[11,2] = "215A_LA PLATA"; $ws3_data[11,3] = 0; $ws3_data[11,4] = "L6173"
$ws3_data[12,0] = "09:21:49"; $ws3_data[12,1] = "10:08"; $ws3_data[12,2] = "215C_LA PLATA"; $ws3_data[12,3] = 47; $ws3_data[12,4] = "L6203"
$ws3_data[13,0] = "08:11:27"; $ws3_data[13,1] = "10:09"; $ws3_data[13,2] = "215C_LA PLATA"; $ws3_data[13,3] = 118; $ws3_data[13,4] = "L6203"
$ws3_data[14,0] = "09:21:49"; $ws3_data[14,1] = "10:22"; $ws3_data[14,2] = "215A_LA PLATA"; $ws3_data[14,3] = 61; $ws3_data[14,4] = "L6173"
$ws3_data[15,0] = "08:29:19"; $ws3_data[15,1] = "10:23"; $ws3_data[15,2] = "215A_LA PLATA"; $ws3_data[15,3] = 114; $ws3_data[15,4] = "L6173"
$ws3_data[16,0] = "09:21:49"; $ws3_data[16,1] = "10:30"; $ws3_data[16,2] = "215B_LP-P MOR-1 Y 57"; $ws3_data[16,3] = 69; $ws3_data[16,4] = "L6173"
$ws3_data[17,0] = "08:39:08"; $ws3_data[17,1] = "10:31"; $ws3_data[17,2] = "215B_LP-P MOR-1 Y 57"; $ws3_data[17,3] = 112; $ws3_data[17,4] = "L6173"
$ws3_data[18,0] = "10:04:17"; $ws3_data[18,1] = "11:44"; $ws3_data[18,2] = "215C_LA PLATA"; $ws3_data[18,3] = 100; $ws3_data[18,4] = "L6203"
$ws3_data[19,0] = "11:53:59"; $ws3_data[19,1] = "12:43"; $ws3_data[19,2] = "215C_LA PLATA"; $ws3_data[19,3] = 50; $ws3_data[19,4] = "L6203"
$ws3_data[20,0] = "10:48:14"; $ws3_data[20,1] = "12:44"; $ws3_data[20,2] = "215C_LA PLATA"; $ws3_data[20,3] = 116; $ws3_data[20,4] = "L6203"
$ws3_data[21,0] = "11:53:59"; $ws3_data[21,1] = "13:08"; $ws3_data[21,2] = "215B_LP-P MOR-1 Y 57"; $ws3_data[21,3] = 75; $ws3_data[21,4] = "L6173"
$ws3_data[22,0] = "11:11:31"; $ws3_data[22,1] = "13:09"; $ws3_data[22,2] = "215B_LP-P MOR-1 Y 57"; $ws3_data[22,3] = 118; $ws3_data[22,4] = "L6173"
$ws3_data[23,0] = "11:53:59"; $ws3_data[23,1] = "13:13"; $ws3_data[23,2] = "215A_LA PLATA"; $ws3_data[23,3] = 80; $ws3_data[23,4] = "L6173"
$ws3_data[24,0] = "11:34:25"; $ws3_data[24,1] = "13:14"; $ws3_data[24,2] = "215A_LA PLATA"; $ws3_data[24,3] = 100; $ws3_data[24,4] = "L6173"
$ws3_data[25,0] = "12:11:45"; $ws3_data[25,1] = "13:53"; $ws3_data[25,2] = "215C_LA PLATA"; $ws3_data[25,3] = 102; $ws3_data[25,4] = "L6203"
$ws3.Range("A6:E31").Value = $ws3_data

